$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2028864896758006
$ws.Range("B1").Value = 0.002539495620456482
$ws.Range("C1").Value = 1.379316785363834
$ws.Range("D1").Value = 0.188939888900502
$ws.Range("E1").Value = 1.570796384046431
$ws.Range("F1").Value = -1.367909888287239
$ws.Range("A2").Value = 0.2484977296296271
$ws.Range("B2").Value = 0.002377709119405816
$ws.Range("C2").Value = 1.380518500463224
$ws.Range("D2").Value = 0.1878999627744357
$ws.Range("E2").Value = 1.570796384423863
$ws.Range("F2").Value = -1.322298649847111
$ws.Range("A3").Value = 0.4528663241188721
$ws.Range("B3").Value = 0.001652798332930176
$ws.Range("C3").Value = 1.385902980825687
$ws.Range("D3").Value = 0.1832404042758736
$ws.Range("E3").Value = 1.570796386115009
$ws.Range("F3").Value = -1.117930062140237
$ws.Range("A4").Value = 0.7403935328885739
$ws.Range("B4").Value = 0.0006329176826605036
$ws.Range("C4").Value = 1.393478433468806
$ws.Range("D4").Value = 0.1766848478679089
$ws.Range("E4").Value = 1.570796388494291
$ws.Range("F4").Value = -0.8304028629126882
$ws.Range("A5").Value = 0.9447621273778194
$ws.Range("B5").Value = -0.00009199310381513776
$ws.Range("C5").Value = 1.398862913831269
$ws.Range("D5").Value = 0.1720252893693467
$ws.Range("E5").Value = 1.570796390185437
$ws.Range("F5").Value = -0.6260342752058142
$ws.Range("A6").Value = 0.9903733673316449
$ws.Range("B6").Value = -0.0002537796048658041
$ws.Range("C6").Value = 1.40006462893066
$ws.Range("D6").Value = 0.1709853632432805
$ws.Range("E6").Value = 1.570796390562869
$ws.Range("F6").Value = -0.5804230367656862
